$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.434937333333333
$ws.Range("H2").Value = 4.304812
$ws.Range("I2").Value = 0.5010808920723563
$ws.Range("J2").Value = 0.5010808920723562
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.23061133333333
$ws.Range("N2").Value = 30.691834
$ws.Range("O2").Value = 0.4855635428718841
$ws.Range("P2").Value = 0.4855635428718841
$ws.Range("Q2").Value = 14.68028614502311
$ws.Range("R2").Value = 132.122575305208
$ws.Range("S2").Value = 0.2433066132200575
$ws.Range("T2").Value = 0.2433066132200574
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.434937333333333
$ws.Range("H3").Value = 4.304812
$ws.Range("I3").Value = 0.5010808920723563
$ws.Range("J3").Value = 0.5010808920723562
$ws.Range("O3").Value = 0.4164864079521221
$ws.Range("P3").Value = 0.4164864079521222
$ws.Range("Q3").Value = 12.59184247665644
$ws.Range("R3").Value = 113.326582289908
$ws.Range("S3").Value = 0.2086933808326606
$ws.Range("T3").Value = 0.2086933808326606
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.434937333333333
$ws.Range("H4").Value = 4.304812
$ws.Range("I4").Value = 0.5010808920723563
$ws.Range("J4").Value = 0.5010808920723562
$ws.Range("M4").Value = 2.034752
$ws.Range("N4").Value = 6.104255999999999
$ws.Range("O4").Value = 0.09657305490303886
$ws.Range("P4").Value = 0.09657305490303887
$ws.Range("Q4").Value = 2.919741608874666
$ws.Range("R4").Value = 26.277674479872
$ws.Range("S4").Value = 0.04839091250096735
$ws.Range("T4").Value = 0.04839091250096735
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.434937333333333
$ws.Range("H5").Value = 4.304812
$ws.Range("I5").Value = 0.5010808920723563
$ws.Range("J5").Value = 0.5010808920723562
$ws.Range("M5").Value = 0.02901266666666667
$ws.Range("N5").Value = 0.087038
$ws.Range("O5").Value = 0.001376994272954919
$ws.Range("P5").Value = 0.001376994272954919
$ws.Range("Q5").Value = 0.04163135853955556
$ws.Range("R5").Value = 0.374682226856
$ws.Range("S5").Value = 0.0006899855186707762
$ws.Range("T5").Value = 0.0006899855186707761
$ws.Range("G6").Value = 0.9964423333333334
$ws.Range("I6").Value = 0.3479582011609289
$ws.Range("J6").Value = 0.3479582011609288
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.23061133333333
$ws.Range("N6").Value = 30.691834
$ws.Range("O6").Value = 0.4855635428718841
$ws.Range("P6").Value = 0.4855635428718841
$ws.Range("Q6").Value = 10.19421422841311
$ws.Range("R6").Value = 91.74792805571801
$ws.Range("S6").Value = 0.1689558169270284
$ws.Range("T6").Value = 0.1689558169270283
$ws.Range("G7").Value = 0.9964423333333334
$ws.Range("I7").Value = 0.3479582011609289
$ws.Range("J7").Value = 0.3479582011609288
$ws.Range("O7").Value = 0.4164864079521221
$ws.Range("P7").Value = 0.4164864079521222
$ws.Range("S7").Value = 0.1449198613189972
$ws.Range("T7").Value = 0.1449198613189972
$ws.Range("G8").Value = 0.9964423333333334
$ws.Range("I8").Value = 0.3479582011609289
$ws.Range("J8").Value = 0.3479582011609288
$ws.Range("M8").Value = 2.034752
$ws.Range("N8").Value = 6.104255999999999
$ws.Range("O8").Value = 0.09657305490303886
$ws.Range("P8").Value = 0.09657305490303887
$ws.Range("Q8").Value = 2.027513030634667
$ws.Range("R8").Value = 18.247617275712
$ws.Range("S8").Value = 0.03360338646467703
$ws.Range("T8").Value = 0.03360338646467702
$ws.Range("G9").Value = 0.9964423333333334
$ws.Range("I9").Value = 0.3479582011609289
$ws.Range("J9").Value = 0.3479582011609288
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.02901266666666667
$ws.Range("N9").Value = 0.087038
$ws.Range("O9").Value = 0.001376994272954919
$ws.Range("P9").Value = 0.001376994272954919
$ws.Range("Q9").Value = 0.02890944926955556
$ws.Range("R9").Value = 0.260185043426
$ws.Range("S9").Value = 0.0004791364502262946
$ws.Range("T9").Value = 0.0004791364502262944
$ws.Range("G10").Value = 0.4323043333333333
$ws.Range("H10").Value = 1.296913
$ws.Range("I10").Value = 0.150960906766715
$ws.Range("J10").Value = 0.1509609067667149
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.23061133333333
$ws.Range("N10").Value = 30.691834
$ws.Range("O10").Value = 0.4855635428718841
$ws.Range("P10").Value = 0.4855635428718841
$ws.Range("Q10").Value = 4.422737612049112
$ws.Range("R10").Value = 39.804638508442
$ws.Range("S10").Value = 0.07330111272479831
$ws.Range("T10").Value = 0.07330111272479828
$ws.Range("G11").Value = 0.4323043333333333
$ws.Range("H11").Value = 1.296913
$ws.Range("I11").Value = 0.150960906766715
$ws.Range("J11").Value = 0.1509609067667149
$ws.Range("O11").Value = 0.4164864079521221
$ws.Range("P11").Value = 0.4164864079521222
$ws.Range("Q11").Value = 3.793551077707444
$ws.Range("R11").Value = 34.141959699367
$ws.Range("S11").Value = 0.06287316580046433
$ws.Range("T11").Value = 0.06287316580046431
$ws.Range("G12").Value = 0.4323043333333333
$ws.Range("H12").Value = 1.296913
$ws.Range("I12").Value = 0.150960906766715
$ws.Range("J12").Value = 0.1509609067667149
$ws.Range("M12").Value = 2.034752
$ws.Range("N12").Value = 6.104255999999999
$ws.Range("O12").Value = 0.09657305490303886
$ws.Range("P12").Value = 0.09657305490303887
$ws.Range("Q12").Value = 0.8796321068586666
$ws.Range("R12").Value = 7.916688961727999
$ws.Range("S12").Value = 0.0145787559373945
$ws.Range("T12").Value = 0.01457875593739449
$ws.Range("G13").Value = 0.4323043333333333
$ws.Range("H13").Value = 1.296913
$ws.Range("I13").Value = 0.150960906766715
$ws.Range("J13").Value = 0.1509609067667149
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.02901266666666667
$ws.Range("N13").Value = 0.087038
$ws.Range("O13").Value = 0.001376994272954919
$ws.Range("P13").Value = 0.001376994272954919
$ws.Range("Q13").Value = 0.01254230152155556
$ws.Range("R13").Value = 0.112880713694
$ws.Range("S13").Value = 0.0002078723040578479
$ws.Range("T13").Value = 0.0002078723040578479
